$wb = $excel.ActiveWorkbook

# Sheet: 展览 (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 37721
$ws1.Range("F5").Value = 780
$ws1.Range("F13").Value = 57
$ws1.Range("F15").Value = 27
$ws1.Range("F16").Value = 660
$ws1.Range("F17").Value = 183
$ws1.Range("F19").Value = 445
$ws1.Range("F21").Value = 94
$ws1.Range("F22").Value = 842
$ws1.Range("F23").Value = 2551
$ws1.Range("F24").Value = 1026
$ws1.Range("F25").Value = 568
$ws1.Range("F26").Value = 110
$ws1.Range("F27").Value = 1167
$ws1.Range("F29").Value = 792
$ws1.Range("F30").Value = 67
$ws1.Range("F31").Value = 1167

# Sheet: 演出 (sheet2)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 413

# Sheet: 本地生活 (sheet3)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 640

# Sheet: 全部类型 (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 640
$ws4.Range("F3").Value = 37721
$ws4.Range("F6").Value = 780
$ws4.Range("F11").Value = 413
$ws4.Range("F19").Value = 57
$ws4.Range("F25").Value = 27
$ws4.Range("F27").Value = 660
$ws4.Range("F28").Value = 183
$ws4.Range("F30").Value = 445
$ws4.Range("F32").Value = 94
$ws4.Range("F33").Value = 842
$ws4.Range("F34").Value = 2551
$ws4.Range("F35").Value = 1026
$ws4.Range("F36").Value = 568
$ws4.Range("F37").Value = 110
$ws4.Range("F38").Value = 1167
$ws4.Range("F41").Value = 792
$ws4.Range("F42").Value = 67
$ws4.Range("F43").Value = 1167
